# Correct the sample weights in column P (rows 2-30)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(11.82, 7.26, 10.02, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 13.06, 0, 12.22, 13.81, 17.77, 14.226, 7.44, 9.96, 0, 9.67, 8.66, 0, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 16).Value = $values[$i]
}

# Hide the helper/raw-data columns D through O (the grouping matches the
# original column-formatting blocks), keep P (and A-C) visible
$ws.Columns("D").Hidden = $true
$ws.Columns("E:I").Hidden = $true
$ws.Columns("J").Hidden = $true
$ws.Columns("K:M").Hidden = $true
$ws.Columns("N").Hidden = $true
$ws.Columns("O").Hidden = $true

# Update the view: no frozen/scrolled topLeftCell, and move the active
# selection away from the P2:P30 block to a single cell T24
$ws.Range("T24").Select() | Out-Null
